$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header columns: SoundEffect (S1) and SEAction (T1)
$ws.Range("S1").Value = "SoundEffect"
$ws.Range("T1").Value = "SEAction"

# Update the selection to match the new active cell in the diff
$ws.Range("U5").Select()
